$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 255.79167
$ws.Range("I53").Value = 136.77777
$ws.Range("J53").Value = 327.2
$ws.Range("K53").Value = 136.77777
$ws.Range("L53").Value = 327.2
$ws.Range("M53").Value = 500.22223
$ws.Range("N53").Value = -1601.2

$ws.Range("H64").Value = 142859700
$ws.Range("I64").Value = 500001540
$ws.Range("J64").Value = 2956
$ws.Range("K64").Value = 500001540
$ws.Range("L64").Value = 2956
$ws.Range("M64").Value = -500001292
$ws.Range("N64").Value = -3452

$ws.Range("H67").Value = 142859700
$ws.Range("I67").Value = 500001540
$ws.Range("J67").Value = 2956
$ws.Range("K67").Value = 500001540
$ws.Range("L67").Value = 2956
$ws.Range("M67").Value = -500001292
$ws.Range("N67").Value = -4672

$ws.Range("H76").Value = 7355.4
$ws.Range("I76").Value = 5234.3335
$ws.Range("K76").Value = 5234.3335
$ws.Range("M76").Value = -4919.3335

$ws.Range("H79").Value = 7355.4
$ws.Range("I79").Value = 5234.3335
$ws.Range("K79").Value = 5234.3335
$ws.Range("M79").Value = -4142.3335

$ws.Range("H129").Value = 5801.2
$ws.Range("I129").Value = 17197.666
$ws.Range("K129").Value = 51592.99800000001
$ws.Range("M129").Value = -46592.99800000001

$ws.Range("H132").Value = 1912.2319
$ws.Range("I132").Value = 1497.2407
$ws.Range("J132").Value = 3406.2
$ws.Range("K132").Value = 4491.7221
$ws.Range("L132").Value = 10218.6
$ws.Range("M132").Value = -1961.7221
$ws.Range("N132").Value = -15278.6

$ws.Range("H137").Value = 5261.294
$ws.Range("I137").Value = 1799.4
$ws.Range("J137").Value = 5858.1724
$ws.Range("K137").Value = 5398.200000000001
$ws.Range("L137").Value = 17574.5172
$ws.Range("M137").Value = -2848.200000000001
$ws.Range("N137").Value = -22674.5172

$ws.Range("H138").Value = 3330.8643
$ws.Range("I138").Value = 2161.818
$ws.Range("J138").Value = 4721.081
$ws.Range("K138").Value = 6485.454000000001
$ws.Range("L138").Value = 14163.243
$ws.Range("M138").Value = -1345.454000000001
$ws.Range("N138").Value = -24443.243

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 295529.9
$ws.Range("I61").Value = 1306.5483
$ws.Range("J61").Value = 3335838
$ws.Range("K61").Value = 1306.5483
$ws.Range("L61").Value = 3335838
$ws.Range("M61").Value = -1094.5483
$ws.Range("N61").Value = -3336262

$ws.Range("H63").Value = 1635
$ws.Range("I63").Value = 1635
$ws.Range("K63").Value = 1635
$ws.Range("M63").Value = -949

$ws.Range("H66").Value = 1635
$ws.Range("I66").Value = 1635
$ws.Range("K66").Value = 8175
$ws.Range("M66").Value = -4743

$ws.Range("H88").Value = 3206.5334
$ws.Range("I88").Value = 3183
$ws.Range("K88").Value = 3183
$ws.Range("M88").Value = -2777

$ws.Range("H91").Value = 3206.5334
$ws.Range("I91").Value = 3183
$ws.Range("K91").Value = 3183
$ws.Range("M91").Value = -1779

$ws.Range("H110").Value = 2217.121
$ws.Range("I110").Value = 2345.476
$ws.Range("J110").Value = 1992.5
$ws.Range("K110").Value = 2345.476
$ws.Range("L110").Value = 1992.5
$ws.Range("M110").Value = -300.4760000000001
$ws.Range("N110").Value = -6082.5

$ws.Range("H136").Value = 295529.9
$ws.Range("I136").Value = 1306.5483
$ws.Range("J136").Value = 3335838
$ws.Range("K136").Value = 3919.6449
$ws.Range("L136").Value = 10007514
$ws.Range("M136").Value = -1369.6449
$ws.Range("N136").Value = -10012614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3362.6897
$ws.Range("I86").Value = 2109.0908
$ws.Range("J86").Value = 4128.778
$ws.Range("K86").Value = 2109.0908
$ws.Range("L86").Value = 4128.778
$ws.Range("M86").Value = -986.0907999999999
$ws.Range("N86").Value = -6374.778

$ws.Range("H89").Value = 3362.6897
$ws.Range("I89").Value = 2109.0908
$ws.Range("J89").Value = 4128.778
$ws.Range("K89").Value = 10545.454
$ws.Range("L89").Value = 20643.89
$ws.Range("M89").Value = -4929.454
$ws.Range("N89").Value = -31875.89

$ws.Range("H134").Value = 970.24
$ws.Range("I134").Value = 686.3333
$ws.Range("J134").Value = 2460.75
$ws.Range("K134").Value = 2058.9999
$ws.Range("L134").Value = 7382.25
$ws.Range("M134").Value = 476.0001000000002
$ws.Range("N134").Value = -12452.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5246.316
$ws.Range("I62").Value = 6200
$ws.Range("J62").Value = 2576
$ws.Range("K62").Value = 6200
$ws.Range("L62").Value = 2576
$ws.Range("M62").Value = -5576
$ws.Range("N62").Value = -3824

$ws.Range("H65").Value = 5246.316
$ws.Range("I65").Value = 6200
$ws.Range("J65").Value = 2576
$ws.Range("K65").Value = 31000
$ws.Range("L65").Value = 12880
$ws.Range("M65").Value = -27880
$ws.Range("N65").Value = -19120

$ws.Range("H122").Value = 923.95
$ws.Range("I122").Value = 620.9091
$ws.Range("J122").Value = 1294.3334
$ws.Range("K122").Value = 1862.7273
$ws.Range("L122").Value = 3883.0002
$ws.Range("M122").Value = 587.2727
$ws.Range("N122").Value = -8783.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 475.18182
$ws.Range("I86").Value = 306.7143
$ws.Range("J86").Value = 770
$ws.Range("K86").Value = 920.1428999999999
$ws.Range("L86").Value = 2310
$ws.Range("M86").Value = 265.8571000000001
$ws.Range("N86").Value = -4682

$ws.Range("H89").Value = 475.18182
$ws.Range("I89").Value = 306.7143
$ws.Range("J89").Value = 770
$ws.Range("K89").Value = 2760.4287
$ws.Range("L89").Value = 6930
$ws.Range("M89").Value = 3167.5713
$ws.Range("N89").Value = -18786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4789.091
$ws.Range("I70").Value = 4633.3335
$ws.Range("J70").Value = 4847.5
$ws.Range("K70").Value = 4633.3335
$ws.Range("L70").Value = 4847.5
$ws.Range("M70").Value = -4363.3335
$ws.Range("N70").Value = -5387.5

$ws.Range("H73").Value = 4789.091
$ws.Range("I73").Value = 4633.3335
$ws.Range("J73").Value = 4847.5
$ws.Range("K73").Value = 4633.3335
$ws.Range("L73").Value = 4847.5
$ws.Range("M73").Value = -3697.3335
$ws.Range("N73").Value = -6719.5

$ws.Range("H80").Value = 2892.9167
$ws.Range("I80").Value = 2203.75
$ws.Range("J80").Value = 3237.5
$ws.Range("K80").Value = 2203.75
$ws.Range("L80").Value = 3237.5
$ws.Range("M80").Value = -1205.75
$ws.Range("N80").Value = -5233.5

$ws.Range("H83").Value = 2892.9167
$ws.Range("I83").Value = 2203.75
$ws.Range("J83").Value = 3237.5
$ws.Range("K83").Value = 11018.75
$ws.Range("L83").Value = 16187.5
$ws.Range("M83").Value = -6026.75
$ws.Range("N83").Value = -26171.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8256.409
$ws.Range("I132").Value = 10499.143
$ws.Range("J132").Value = 4331.625
$ws.Range("K132").Value = 31497.429
$ws.Range("L132").Value = 12994.875
$ws.Range("M132").Value = -28967.429
$ws.Range("N132").Value = -18054.875
